# Apply crypto price/volume updates as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "68.605.18"
$ws.Range("E2").Value = "  -0.88%  "

$ws.Range("D3").Value = "3.910.98"
$ws.Range("E3").Value = "  +2.58%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'602.96"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D6").Value = "'166.50"
$ws.Range("E6").Value = "  +1.29%  "

$ws.Range("D7").Value = "3.910.41"
$ws.Range("E7").Value = "  +2.56%  "

$ws.Range("E9").Value = "  -1.42%  "

$ws.Range("D10").Value = "'0.168"
$ws.Range("E10").Value = "  -1.88%  "

$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("D12").Value = "'0.461"
$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("E13").Value = "  +3.28%  "

$ws.Range("D14").Value = "'37.41"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").Value = "4.562.91"
$ws.Range("E15").Value = "  +2.48%  "

$ws.Range("D16").Value = "3.905.08"
$ws.Range("E16").Value = "  +2.65%  "

$ws.Range("D17").Value = "68.731.90"
$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").Value = "'17.30"
$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("E20").Value = "  -2.40%  "

$ws.Range("D21").Value = "'11.07"
$ws.Range("E21").Value = "  -2.27%  "

$ws.Range("D22").Value = "'488.16"
$ws.Range("E22").Value = "  -0.40%  "

$ws.Range("D23").Value = "'0.728"
$ws.Range("E23").Value = "  +0.52%  "

$ws.Range("D24").Value = "'0.0000166"
$ws.Range("E24").Value = "  +6.48%  "

$ws.Range("D25").Value = "'84.70"
$ws.Range("E25").Value = "  -0.23%  "

$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "  -2.02%  "

$ws.Range("D27").Value = "'12.08"
$ws.Range("E27").Value = "  -1.74%  "

$ws.Range("D28").Value = "'10.14"
$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  -1.07%  "

$ws.Range("D31").Value = "4.061.49"
$ws.Range("E31").Value = "  +2.65%  "

$ws.Range("E32").Value = "  -1.22%  "

$ws.Range("E33").Value = "  -3.94%  "

$ws.Range("D34").Value = "'31.93"
$ws.Range("E34").Value = "  -0.34%  "

$ws.Range("D35").Value = "3.864.04"
$ws.Range("E35").Value = "  +2.83%  "

$ws.Range("E36").Value = "  -0.46%  "

$ws.Range("E37").Value = "  +1.70%  "

$ws.Range("D38").Value = "'5.94"
$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  -1.12%  "

$ws.Range("D40").Value = "'3.19"
$ws.Range("E40").Value = "  +5.02%  "

$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").Value = "'0.316"
$ws.Range("E42").Value = "  -1.97%  "

$ws.Range("D43").Value = "'431.21"
$ws.Range("E43").Value = "  +1.59%  "

$ws.Range("D44").Value = "'1.99"
$ws.Range("E44").Value = "  -0.75%  "

$ws.Range("D45").Value = "'48.14"
$ws.Range("E45").Value = "  -1.07%  "

$ws.Range("D46").Value = "'8.54"
$ws.Range("E46").Value = "  +1.66%  "

$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").Value = "'142.52"
$ws.Range("E48").Value = "  +0.78%  "

$ws.Range("D49").Value = "2.811.78"
$ws.Range("E49").Value = "  -0.77%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000266"
$ws.Range("E50").Value = "  +17.43%  "

$ws.Range("D51").Value = "'0.0354"
$ws.Range("E51").Value = "  +0.81%  "
